$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 474.1111
$ws.Range("I9").Value = 389.5
$ws.Range("J9").Value = 643.3333
$ws.Range("K9").Value = 389.5
$ws.Range("L9").Value = 643.3333
$ws.Range("M9").Value = -220.5
$ws.Range("N9").Value = -981.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3932.2683
$ws.Range("I45").Value = 2727.1052
$ws.Range("K45").Value = 2727.1052
$ws.Range("M45").Value = -2350.1052
$ws.Range("H61").Value = 37040890
$ws.Range("I61").Value = 3194.9524
$ws.Range("K61").Value = 3194.9524
$ws.Range("M61").Value = -2982.9524
$ws.Range("H110").Value = 22223606
$ws.Range("I110").Value = 1465.1538
$ws.Range("K110").Value = 1465.1538
$ws.Range("M110").Value = 579.8462
$ws.Range("H136").Value = 37040890
$ws.Range("I136").Value = 3194.9524
$ws.Range("K136").Value = 9584.8572
$ws.Range("M136").Value = -7034.8572

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 167350.17
$ws.Range("I86").Value = 333700.34
$ws.Range("K86").Value = 333700.34
$ws.Range("M86").Value = -332577.34
$ws.Range("H89").Value = 167350.17
$ws.Range("I89").Value = 333700.34
$ws.Range("K89").Value = 1668501.7
$ws.Range("M89").Value = -1662885.7
$ws.Range("H134").Value = 5957020.5
$ws.Range("I134").Value = 9618270
$ws.Range("K134").Value = 28854810
$ws.Range("M134").Value = -28852275

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5148.9287
$ws.Range("I16").Value = 881.4167
$ws.Range("J16").Value = 8349.5625
$ws.Range("K16").Value = 881.4167
$ws.Range("L16").Value = 8349.5625
$ws.Range("M16").Value = -594.4167
$ws.Range("N16").Value = -8923.5625
$ws.Range("H31").Value = 4507.8184
$ws.Range("I31").Value = 1641.3846
$ws.Range("K31").Value = 1641.3846
$ws.Range("M31").Value = -1346.3846
$ws.Range("H34").Value = 4507.8184
$ws.Range("I34").Value = 1641.3846
$ws.Range("K34").Value = 1641.3846
$ws.Range("M34").Value = -1439.3846
$ws.Range("H58").Value = 4798.575
$ws.Range("J58").Value = 6480.4
$ws.Range("L58").Value = 6480.4
$ws.Range("N58").Value = -6886.4
$ws.Range("H69").Value = 20706.23
$ws.Range("I69").Value = 16652.818
$ws.Range("J69").Value = 43000
$ws.Range("K69").Value = 16652.818
$ws.Range("L69").Value = 43000
$ws.Range("M69").Value = -15903.818
$ws.Range("N69").Value = -44498
$ws.Range("H72").Value = 20706.23
$ws.Range("I72").Value = 16652.818
$ws.Range("J72").Value = 43000
$ws.Range("K72").Value = 49958.454
$ws.Range("L72").Value = 129000
$ws.Range("M72").Value = -46214.454
$ws.Range("N72").Value = -136488
$ws.Range("H113").Value = 5148.9287
$ws.Range("I113").Value = 881.4167
$ws.Range("J113").Value = 8349.5625
$ws.Range("K113").Value = 881.4167
$ws.Range("L113").Value = 8349.5625
$ws.Range("M113").Value = 1288.5833
$ws.Range("N113").Value = -12689.5625
$ws.Range("H122").Value = 1707.05
$ws.Range("I122").Value = 1092.6875
$ws.Range("J122").Value = 4164.5
$ws.Range("K122").Value = 3278.0625
$ws.Range("L122").Value = 12493.5
$ws.Range("M122").Value = -828.0625
$ws.Range("N122").Value = -17393.5
$ws.Range("H132").Value = 4658.7837
$ws.Range("I132").Value = 3489.524
$ws.Range("J132").Value = 6193.4375
$ws.Range("K132").Value = 10468.572
$ws.Range("L132").Value = 18580.3125
$ws.Range("M132").Value = -7938.572
$ws.Range("N132").Value = -23640.3125
$ws.Range("H134").Value = 4550.8696
$ws.Range("I134").Value = 2608.25
$ws.Range("K134").Value = 7824.75
$ws.Range("M134").Value = -5289.75
$ws.Range("H136").Value = 4798.575
$ws.Range("J136").Value = 6480.4
$ws.Range("L136").Value = 19441.2
$ws.Range("N136").Value = -24541.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 104486.54
$ws.Range("J2").Value = 275715.88
$ws.Range("L2").Value = 1654295.28
$ws.Range("N2").Value = -1654521.28
$ws.Range("H39").Value = 13814.167
$ws.Range("J39").Value = 13777
$ws.Range("L39").Value = 41331
$ws.Range("N39").Value = -41919
$ws.Range("H55").Value = 5008360
$ws.Range("J55").Value = 5271802.5
$ws.Range("L55").Value = 15815407.5
$ws.Range("N55").Value = -15815761.5
$ws.Range("H97").Value = 502.83334
$ws.Range("I97").Value = 662.6667
$ws.Range("K97").Value = 1988.0001
$ws.Range("M97").Value = -1492.0001
$ws.Range("H98").Value = 1488.1
$ws.Range("I98").Value = 567.3333
$ws.Range("J98").Value = 1882.7142
$ws.Range("K98").Value = 1701.9999
$ws.Range("L98").Value = 5648.142599999999
$ws.Range("M98").Value = -203.9999
$ws.Range("N98").Value = -8644.142599999999
$ws.Range("H140").Value = 64218.5
$ws.Range("I140").Value = 92136
$ws.Range("K140").Value = 276408
$ws.Range("M140").Value = -271228

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 38247.5
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 38247.5
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 38247.5
$ws.Range("N15").Value = -38823.5
$ws.Range("M15").Value = $null
$ws.Range("H33").Value = 47976
$ws.Range("J33").Value = 47976
$ws.Range("L33").Value = 47976
$ws.Range("N33").Value = -48480
$ws.Range("H81").Value = 38247.5
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 38247.5
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 38247.5
$ws.Range("N81").Value = -40243.5
$ws.Range("M81").Value = $null
$ws.Range("H84").Value = 38247.5
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 38247.5
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 114742.5
$ws.Range("N84").Value = -124726.5
$ws.Range("M84").Value = $null
$ws.Range("H102").Value = 6217.9375
$ws.Range("I102").Value = 6104.2856
$ws.Range("K102").Value = 6104.2856
$ws.Range("M102").Value = -4482.2856

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3698.0881
$ws.Range("I40").Value = 2776.3215
$ws.Range("J40").Value = 7999.6665
$ws.Range("K40").Value = 2776.3215
$ws.Range("L40").Value = 7999.6665
$ws.Range("M40").Value = -2640.3215
$ws.Range("N40").Value = -8271.666499999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 649.26086
$ws.Range("I100").Value = 422.89474
$ws.Range("J100").Value = 1724.5
$ws.Range("K100").Value = 845.78948
$ws.Range("L100").Value = 3449
$ws.Range("M100").Value = -304.78948
$ws.Range("N100").Value = -4531
$ws.Range("H113").Value = 1266.2812
$ws.Range("I113").Value = 1094.3889
$ws.Range("J113").Value = 1487.2858
$ws.Range("K113").Value = 3283.1667
$ws.Range("L113").Value = 4461.857400000001
$ws.Range("M113").Value = -1113.1667
$ws.Range("N113").Value = -8801.857400000001
$ws.Range("H122").Value = 193704.67
$ws.Range("I122").Value = 237401.06
$ws.Range("K122").Value = 712203.1799999999
$ws.Range("M122").Value = -709753.1799999999
$ws.Range("H126").Value = 743.26666
$ws.Range("J126").Value = 966.5
$ws.Range("L126").Value = 2899.5
$ws.Range("N126").Value = -7839.5

Write-Output "Applied all Sargatanas_Profits updates"